$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Cells.Item(2, 9).Value = "b"
$ws.Cells.Item(2, 10).Value = "Acknowledge (Backchannel)"
$ws.Cells.Item(4, 9).Value = "sd"
$ws.Cells.Item(4, 10).Value = "Statement-non-opinion"
$ws.Cells.Item(31, 9).Value = "aa"
$ws.Cells.Item(31, 10).Value = "Agree/Accept"
$ws.Cells.Item(36, 9).Value = "b"
$ws.Cells.Item(36, 10).Value = "Acknowledge (Backchannel)"
$ws.Cells.Item(59, 9).Value = "sd"
$ws.Cells.Item(59, 10).Value = "Statement-non-opinion"
$ws.Cells.Item(65, 9).Value = "sd"
$ws.Cells.Item(65, 10).Value = "Statement-non-opinion"
$ws.Cells.Item(68, 9).Value = "sv"
$ws.Cells.Item(68, 10).Value = "Statement-opinion"
$ws.Cells.Item(71, 9).Value = "sv"
$ws.Cells.Item(71, 10).Value = "Statement-opinion"
$ws.Cells.Item(74, 9).Value = "aa"
$ws.Cells.Item(74, 10).Value = "Agree/Accept"
$ws.Cells.Item(76, 9).Value = "sv"
$ws.Cells.Item(76, 10).Value = "Statement-opinion"
$ws.Cells.Item(80, 9).Value = "b"
$ws.Cells.Item(80, 10).Value = "Acknowledge (Backchannel)"
$ws.Cells.Item(92, 9).Value = "aa"
$ws.Cells.Item(92, 10).Value = "Agree/Accept"
$ws.Cells.Item(97, 9).Value = "b"
$ws.Cells.Item(97, 10).Value = "Acknowledge (Backchannel)"
$ws.Cells.Item(101, 9).Value = "sv"
$ws.Cells.Item(101, 10).Value = "Statement-opinion"
$ws.Cells.Item(106, 9).Value = "sd"
$ws.Cells.Item(106, 10).Value = "Statement-non-opinion"
$ws.Cells.Item(120, 9).Value = "sd"
$ws.Cells.Item(120, 10).Value = "Statement-non-opinion"
$ws.Cells.Item(125, 9).Value = "sd"
$ws.Cells.Item(125, 10).Value = "Statement-non-opinion"
$ws.Cells.Item(128, 9).Value = "sd"
$ws.Cells.Item(128, 10).Value = "Statement-non-opinion"
$ws.Cells.Item(132, 9).Value = "b"
$ws.Cells.Item(132, 10).Value = "Acknowledge (Backchannel)"
$ws.Cells.Item(136, 9).Value = "aa"
$ws.Cells.Item(136, 10).Value = "Agree/Accept"
$ws.Cells.Item(139, 9).Value = "sd"
$ws.Cells.Item(139, 10).Value = "Statement-non-opinion"
$ws.Cells.Item(142, 9).Value = "%"
$ws.Cells.Item(142, 10).Value = "Uninterpretable"
$ws.Cells.Item(157, 9).Value = "sd"
$ws.Cells.Item(157, 10).Value = "Statement-non-opinion"
$ws.Cells.Item(166, 9).Value = "sd"
$ws.Cells.Item(166, 10).Value = "Statement-non-opinion"
$ws.Cells.Item(168, 9).Value = "aa"
$ws.Cells.Item(168, 10).Value = "Agree/Accept"
$ws.Cells.Item(171, 9).Value = "aa"
$ws.Cells.Item(171, 10).Value = "Agree/Accept"
$ws.Cells.Item(181, 9).Value = "b"
$ws.Cells.Item(181, 10).Value = "Acknowledge (Backchannel)"
$ws.Cells.Item(182, 9).Value = "b"
$ws.Cells.Item(182, 10).Value = "Acknowledge (Backchannel)"
$ws.Cells.Item(183, 9).Value = "aa"
$ws.Cells.Item(183, 10).Value = "Agree/Accept"
$ws.Cells.Item(189, 9).Value = "b"
$ws.Cells.Item(189, 10).Value = "Acknowledge (Backchannel)"
$ws.Cells.Item(191, 9).Value = "sd"
$ws.Cells.Item(191, 10).Value = "Statement-non-opinion"
$ws.Cells.Item(194, 9).Value = "sd"
$ws.Cells.Item(194, 10).Value = "Statement-non-opinion"
$ws.Cells.Item(214, 9).Value = "sd"
$ws.Cells.Item(214, 10).Value = "Statement-non-opinion"
$ws.Cells.Item(228, 9).Value = "sd"
$ws.Cells.Item(228, 10).Value = "Statement-non-opinion"
$ws.Cells.Item(230, 9).Value = "aa"
$ws.Cells.Item(230, 10).Value = "Agree/Accept"
$ws.Cells.Item(251, 9).Value = "sd"
$ws.Cells.Item(251, 10).Value = "Statement-non-opinion"
$ws.Cells.Item(256, 9).Value = "sd"
$ws.Cells.Item(256, 10).Value = "Statement-non-opinion"
$ws.Cells.Item(257, 9).Value = "sd"
$ws.Cells.Item(257, 10).Value = "Statement-non-opinion"
$ws.Cells.Item(267, 9).Value = "sd"
$ws.Cells.Item(267, 10).Value = "Statement-non-opinion"
$ws.Cells.Item(268, 9).Value = "sd"
$ws.Cells.Item(268, 10).Value = "Statement-non-opinion"
$ws.Cells.Item(275, 9).Value = "%"
$ws.Cells.Item(275, 10).Value = "Uninterpretable"
$ws.Cells.Item(278, 9).Value = "sv"
$ws.Cells.Item(278, 10).Value = "Statement-opinion"
$ws.Cells.Item(282, 9).Value = "aa"
$ws.Cells.Item(282, 10).Value = "Agree/Accept"
$ws.Cells.Item(286, 9).Value = "%"
$ws.Cells.Item(286, 10).Value = "Uninterpretable"
$ws.Cells.Item(287, 9).Value = "b"
$ws.Cells.Item(287, 10).Value = "Acknowledge (Backchannel)"
$ws.Cells.Item(291, 9).Value = "ba"
$ws.Cells.Item(291, 10).Value = "Appreciation"
$ws.Cells.Item(295, 9).Value = "sv"
$ws.Cells.Item(295, 10).Value = "Statement-opinion"
$ws.Cells.Item(296, 9).Value = "aa"
$ws.Cells.Item(296, 10).Value = "Agree/Accept"
$ws.Cells.Item(303, 9).Value = "sv"
$ws.Cells.Item(303, 10).Value = "Statement-opinion"
$ws.Cells.Item(306, 9).Value = "sd"
$ws.Cells.Item(306, 10).Value = "Statement-non-opinion"
$ws.Cells.Item(309, 9).Value = "sv"
$ws.Cells.Item(309, 10).Value = "Statement-opinion"
$ws.Cells.Item(338, 9).Value = "aa"
$ws.Cells.Item(338, 10).Value = "Agree/Accept"
$ws.Cells.Item(342, 9).Value = "sd"
$ws.Cells.Item(342, 10).Value = "Statement-non-opinion"
$ws.Cells.Item(347, 9).Value = "sd"
$ws.Cells.Item(347, 10).Value = "Statement-non-opinion"
$ws.Cells.Item(348, 9).Value = "sd"
$ws.Cells.Item(348, 10).Value = "Statement-non-opinion"
$ws.Cells.Item(352, 9).Value = "sd"
$ws.Cells.Item(352, 10).Value = "Statement-non-opinion"
$ws.Cells.Item(353, 9).Value = "aa"
$ws.Cells.Item(353, 10).Value = "Agree/Accept"
$ws.Cells.Item(354, 9).Value = "sd"
$ws.Cells.Item(354, 10).Value = "Statement-non-opinion"
$ws.Cells.Item(355, 9).Value = "aa"
$ws.Cells.Item(355, 10).Value = "Agree/Accept"
$ws.Cells.Item(363, 9).Value = "b"
$ws.Cells.Item(363, 10).Value = "Acknowledge (Backchannel)"
$ws.Cells.Item(387, 9).Value = "aa"
$ws.Cells.Item(387, 10).Value = "Agree/Accept"
$ws.Cells.Item(388, 9).Value = "sd"
$ws.Cells.Item(388, 10).Value = "Statement-non-opinion"
$ws.Cells.Item(390, 9).Value = "%"
$ws.Cells.Item(390, 10).Value = "Uninterpretable"
$ws.Cells.Item(392, 9).Value = "b"
$ws.Cells.Item(392, 10).Value = "Acknowledge (Backchannel)"
$ws.Cells.Item(393, 9).Value = "aa"
$ws.Cells.Item(393, 10).Value = "Agree/Accept"
$ws.Cells.Item(401, 9).Value = "sd"
$ws.Cells.Item(401, 10).Value = "Statement-non-opinion"
$ws.Cells.Item(403, 9).Value = "sd"
$ws.Cells.Item(403, 10).Value = "Statement-non-opinion"
$ws.Cells.Item(405, 9).Value = "sd"
$ws.Cells.Item(405, 10).Value = "Statement-non-opinion"
$ws.Cells.Item(412, 9).Value = "sd"
$ws.Cells.Item(412, 10).Value = "Statement-non-opinion"
$ws.Cells.Item(415, 9).Value = "sd"
$ws.Cells.Item(415, 10).Value = "Statement-non-opinion"
$ws.Cells.Item(423, 9).Value = "aa"
$ws.Cells.Item(423, 10).Value = "Agree/Accept"
$ws.Cells.Item(449, 9).Value = "aa"
$ws.Cells.Item(449, 10).Value = "Agree/Accept"
$ws.Cells.Item(471, 9).Value = "ba"
$ws.Cells.Item(471, 10).Value = "Appreciation"
$ws.Cells.Item(473, 9).Value = "sd"
$ws.Cells.Item(473, 10).Value = "Statement-non-opinion"
